# Commit message: "Attempted adding loading screen in download"
#
# Net effect observed in the diff: the "01:00:00" shared string entry is
# dropped, which corresponds to clearing the "OFFICIAL BUSINESS DEPARTURE"
# (column K) and "OFFICIAL BUSINESS ARRIVAL" (column N) cells for the
# attendance rows 8 through 18 -- those were the only cells holding that
# value. The neighboring "OFFICIAL BUSINESS TIME START"/"TIME END" columns
# (L and M) keep their original values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K8:K18").ClearContents()
$ws.Range("N8:N18").ClearContents()
